# Update of league bases (Denmark Division 3) - 28-05-2024 20:56
# The underlying source data re-fetched a handful of fixtures whose rows
# ended up re-ordered in the new export. Re-apply the same row content
# shuffle to the existing rows (columns B:AD only; column A, the row
# index, is untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) contents of every row that participates
# in a re-shuffle, BEFORE any writes happen, since several groups are
# cyclic permutations (not simple pairwise swaps).
$rows = 32,33,95,96,111,112,113,114,115,116,117,124,125,141,142,148,149,150
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = $ws.Range("B$($r):AD$($r)").Value()
}

# Destination row -> source row (source row's ORIGINAL content is copied
# into destination row).
$mapping = @{
    32  = 33
    33  = 32
    95  = 96
    96  = 95
    111 = 112
    112 = 113
    113 = 114
    114 = 115
    115 = 111
    116 = 117
    117 = 116
    124 = 125
    125 = 124
    141 = 142
    142 = 141
    148 = 149
    149 = 150
    150 = 148
}

foreach ($dest in $rows) {
    $src = $mapping[$dest]
    $ws.Range("B$($dest):AD$($dest)").Value = $snapshot[$src]
}
